$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A, G, H to match the new layout.
# Excel's ColumnWidth property (character units) gets stored in the XML
# "width" attribute with a constant +5/6 padding added for this font, so
# subtract that offset here to land on the exact target stored widths
# (21, 15, 22).
$colPad = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 21 - $colPad
$ws.Columns.Item(7).ColumnWidth = 15 - $colPad
$ws.Columns.Item(8).ColumnWidth = 22 - $colPad

# New data rows 8-12 (added at the bottom of the sheet)
$data = @(
    @("address {{address}}", "f4 address", " ", "first f4", "last f4", " ", "{{pin}}", "{{in_aadhaar}}"),
    @("f1 adrress", "01 city", " ", "f1 first", "f1 last", "{{phone}}", "{{date_time}}", "{{date_time}}"),
    @("ddaddress", "dd city", " ", "ddirst", "{{address}}", " ", " ", "{{address}}"),
    @("address", "HCM city", " ", "first", "last", "039494944", "2", "{{credit_debit_cvv}}"),
    @("address", "H C M City ", " ", "firsst", "last", " ", "222", " ")
)

# Cells whose text looks like a plain number and would otherwise be
# auto-converted to a numeric value by Excel; force them to stay text
# (preserving e.g. the leading zero in "039494944").
$textCells = @{
    "11,6" = $true
    "11,7" = $true
    "12,7" = $true
}

$rowIndex = 8
foreach ($row in $data) {
    for ($col = 1; $col -le 8; $col++) {
        $cell = $ws.Cells.Item($rowIndex, $col)
        $value = $row[$col - 1]
        if ($textCells.ContainsKey("$rowIndex,$col")) {
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
    $rowIndex++
}
